$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.473.09'
$ws.Range('E2').Value = '  -0.60%  '
$ws.Range('D3').Value = '3.418.66'
$ws.Range('E3').Value = '  -1.81%  '
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '0.998'
$c.Style = 'Normal'
$ws.Range('E4').Value = '  -0.26%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '576.70'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -0.02%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '164.23'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +1.87%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = '3.409.98'
$ws.Range('E8').Value = '  -1.99%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.552'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  -4.49%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '7.30'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  +0.88%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.120'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  -2.33%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.422'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  -3.98%  '
$ws.Range('D13').Value = '3.997.36'
$ws.Range('E13').Value = '  -2.01%  '
$ws.Range('E14').Value = '  +0.27%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '26.95'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  -2.57%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '0.0000173'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  -1.93%  '
$ws.Range('D17').Value = '64.414.04'
$ws.Range('E17').Value = '  -0.86%  '
$ws.Range('D18').Value = '3.447.32'
$ws.Range('E18').Value = '  -3.37%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '6.17'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  -1.14%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '13.52'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -2.65%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '375.42'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -1.90%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '7.83'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -2.17%  '
$ws.Range('E23').Value = '  -0.05%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '70.43'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -3.25%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '0.516'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -3.39%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '0.0000116'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -4.74%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '9.52'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -3.94%  '
$ws.Range('E28').Value = '  -0.68%  '
$ws.Range('E29').Value = '  -0.23%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '6.10'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -0.74%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '1.40'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  -3.44%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '2.01'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -0.55%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +0.09%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '22.93'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -2.29%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '7.04'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  -0.76%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '1.48'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -7.00%  '
$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '159.54'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -1.24%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.861'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  +5.60%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '1.83'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -2.96%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.0721'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -4.58%  '
$ws.Range('B41').Value = 'EnergySwap'
$ws.Range('C41').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '25.77'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -4.86%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '2.770.08'
$ws.Range('E42').Value = '  -3.93%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '42.67'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -0.70%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '6.45'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -1.58%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '25.84'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -0.85%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '4.37'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -3.75%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '0.0304'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -1.94%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '2.45'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -0.13%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '328.73'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -0.97%  '
$ws.Range('E50').Value = '  -2.86%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '6.29'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -3.04%  '
